# The "Nio" column (S) was a spurious/duplicate column and is removed
# entirely. Deleting the whole column shifts T:AB left by one (into S:AA),
# which is exactly what the diff shows (dimension A1:AB9 -> A1:AA9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S1").EntireColumn.Delete()

# Recalculated descriptive statistics for the "Ni" column (M) after the
# "Nio" values were folded into it.
$ws.Range("M3").Value = 5.594993006993007
$ws.Range("M4").Value = 5.893864266583386
$ws.Range("M7").Value = 5.3
$ws.Range("M8").Value = 10
